$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "1497"
$ws.Range("C4").ClearFormats()
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "11085651.35"
$ws.Range("D4").ClearFormats()

$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "1054"
$ws.Range("C6").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6900197.66"
$ws.Range("D6").ClearFormats()

$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "315"
$ws.Range("C9").ClearFormats()
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1039706.10"
$ws.Range("D9").ClearFormats()

$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "531"
$ws.Range("C11").ClearFormats()
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3815458.66"
$ws.Range("D11").ClearFormats()

$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "262"
$ws.Range("C12").ClearFormats()
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1799569.34"
$ws.Range("D12").ClearFormats()

$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "531"
$ws.Range("C44").ClearFormats()
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2026889.17"
$ws.Range("D44").ClearFormats()

$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "205"
$ws.Range("C45").ClearFormats()
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1288625.38"
$ws.Range("D45").ClearFormats()

$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "289"
$ws.Range("C46").ClearFormats()
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1345741.02"
$ws.Range("D46").ClearFormats()

$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "1150"
$ws.Range("C51").ClearFormats()
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8368747.48"
$ws.Range("D51").ClearFormats()

$ws.Range("C52").NumberFormat = "@"
$ws.Range("C52").Value = "797"
$ws.Range("C52").ClearFormats()
$ws.Range("D52").NumberFormat = "@"
$ws.Range("D52").Value = "5102970.45"
$ws.Range("D52").ClearFormats()

$ws.Range("C75").NumberFormat = "@"
$ws.Range("C75").Value = "300"
$ws.Range("C75").ClearFormats()
$ws.Range("D75").NumberFormat = "@"
$ws.Range("D75").Value = "1053695.80"
$ws.Range("D75").ClearFormats()

$ws.Range("C76").NumberFormat = "@"
$ws.Range("C76").Value = "528"
$ws.Range("C76").ClearFormats()
$ws.Range("D76").NumberFormat = "@"
$ws.Range("D76").Value = "3650765.72"
$ws.Range("D76").ClearFormats()

$ws.Range("C77").NumberFormat = "@"
$ws.Range("C77").Value = "314"
$ws.Range("C77").ClearFormats()
$ws.Range("D77").NumberFormat = "@"
$ws.Range("D77").Value = "2595570.99"
$ws.Range("D77").ClearFormats()

$ws.Range("C80").NumberFormat = "@"
$ws.Range("C80").Value = "453"
$ws.Range("C80").ClearFormats()
$ws.Range("D80").NumberFormat = "@"
$ws.Range("D80").Value = "1450372.96"
$ws.Range("D80").ClearFormats()

$ws.Range("C82").NumberFormat = "@"
$ws.Range("C82").Value = "1272"
$ws.Range("C82").ClearFormats()
$ws.Range("D82").NumberFormat = "@"
$ws.Range("D82").Value = "10051714.73"
$ws.Range("D82").ClearFormats()

$ws.Range("C83").NumberFormat = "@"
$ws.Range("C83").Value = "668"
$ws.Range("C83").ClearFormats()
$ws.Range("D83").NumberFormat = "@"
$ws.Range("D83").Value = "4526673.10"
$ws.Range("D83").ClearFormats()

$ws.Range("C86").NumberFormat = "@"
$ws.Range("C86").Value = "997"
$ws.Range("C86").ClearFormats()
$ws.Range("D86").NumberFormat = "@"
$ws.Range("D86").Value = "2839424.58"
$ws.Range("D86").ClearFormats()

$ws.Range("C89").NumberFormat = "@"
$ws.Range("C89").Value = "1397"
$ws.Range("C89").ClearFormats()
$ws.Range("D89").NumberFormat = "@"
$ws.Range("D89").Value = "9423818.91"
$ws.Range("D89").ClearFormats()

$ws.Range("C91").NumberFormat = "@"
$ws.Range("C91").Value = "994"
$ws.Range("C91").ClearFormats()
$ws.Range("D91").NumberFormat = "@"
$ws.Range("D91").Value = "6109118.94"
$ws.Range("D91").ClearFormats()

$ws.Range("C104").NumberFormat = "@"
$ws.Range("C104").Value = "1694"
$ws.Range("C104").ClearFormats()
$ws.Range("D104").NumberFormat = "@"
$ws.Range("D104").Value = "9464993.02"
$ws.Range("D104").ClearFormats()

$ws.Range("C106").NumberFormat = "@"
$ws.Range("C106").Value = "1653"
$ws.Range("C106").ClearFormats()
$ws.Range("D106").NumberFormat = "@"
$ws.Range("D106").Value = "8665387.85"
$ws.Range("D106").ClearFormats()
